$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "60-19="
$t.Cell(1,2).Range.Text = "4+30="
$t.Cell(1,3).Range.Text = "26+59="
$t.Cell(1,4).Range.Text = "62+30="
$t.Cell(1,5).Range.Text = "31+18="
$t.Cell(2,1).Range.Text = "6+82="
$t.Cell(2,2).Range.Text = "88+0="
$t.Cell(2,3).Range.Text = "23+4="
$t.Cell(2,4).Range.Text = "16+59="
$t.Cell(2,5).Range.Text = "93-6="
$t.Cell(3,1).Range.Text = "71-52="
$t.Cell(3,2).Range.Text = "16+12="
$t.Cell(3,3).Range.Text = "16-10="
$t.Cell(3,4).Range.Text = "82-24="
$t.Cell(3,5).Range.Text = "80-27="
$t.Cell(4,1).Range.Text = "5+16="
$t.Cell(4,2).Range.Text = "1+44="
$t.Cell(4,3).Range.Text = "52-22="
$t.Cell(4,4).Range.Text = "30+1="
$t.Cell(4,5).Range.Text = "75+13="
$t.Cell(5,1).Range.Text = "83-37="
$t.Cell(5,2).Range.Text = "61+35="
$t.Cell(5,3).Range.Text = "50-22="
$t.Cell(5,4).Range.Text = "9+28="
$t.Cell(5,5).Range.Text = "62-7="
$t.Cell(6,1).Range.Text = "56-48="
$t.Cell(6,2).Range.Text = "70+21="
$t.Cell(6,3).Range.Text = "51+20="
$t.Cell(6,4).Range.Text = "15+73="
$t.Cell(6,5).Range.Text = "90-68="
$t.Cell(7,1).Range.Text = "72-64="
$t.Cell(7,2).Range.Text = "92-0="
$t.Cell(7,3).Range.Text = "58-8="
$t.Cell(7,4).Range.Text = "82-81="
$t.Cell(7,5).Range.Text = "77+8="
$t.Cell(8,1).Range.Text = "90-9="
$t.Cell(8,2).Range.Text = "10+68="
$t.Cell(8,3).Range.Text = "21+1="
$t.Cell(8,4).Range.Text = "18+47="
$t.Cell(8,5).Range.Text = "54+5="
$t.Cell(9,1).Range.Text = "5+70="
$t.Cell(9,2).Range.Text = "97-27="
$t.Cell(9,3).Range.Text = "96-41="
$t.Cell(9,4).Range.Text = "15+81="
$t.Cell(9,5).Range.Text = "80-67="
$t.Cell(10,1).Range.Text = "0+76="
$t.Cell(10,2).Range.Text = "74-19="
$t.Cell(10,3).Range.Text = "5+86="
$t.Cell(10,4).Range.Text = "10+34="
$t.Cell(10,5).Range.Text = "39+25="
$t.Cell(11,1).Range.Text = "11+7="
$t.Cell(11,2).Range.Text = "15+57="
$t.Cell(11,3).Range.Text = "54-21="
$t.Cell(11,4).Range.Text = "85-69="
$t.Cell(11,5).Range.Text = "30-24="
$t.Cell(12,1).Range.Text = "42+9="
$t.Cell(12,2).Range.Text = "37+24="
$t.Cell(12,3).Range.Text = "56-31="
$t.Cell(12,4).Range.Text = "86-35="
$t.Cell(12,5).Range.Text = "80-65="
$t.Cell(13,1).Range.Text = "88-71="
$t.Cell(13,2).Range.Text = "6+68="
$t.Cell(13,3).Range.Text = "69-25="
$t.Cell(13,4).Range.Text = "85+10="
$t.Cell(13,5).Range.Text = "11+58="
$t.Cell(14,1).Range.Text = "97-75="
$t.Cell(14,2).Range.Text = "65+33="
$t.Cell(14,3).Range.Text = "26-24="
$t.Cell(14,4).Range.Text = "25+17="
$t.Cell(14,5).Range.Text = "86-31="
$t.Cell(15,1).Range.Text = "14+66="
$t.Cell(15,2).Range.Text = "67+12="
$t.Cell(15,3).Range.Text = "75-47="
$t.Cell(15,4).Range.Text = "21+59="
$t.Cell(15,5).Range.Text = "26+51="
$t.Cell(16,1).Range.Text = "28+33="
$t.Cell(16,2).Range.Text = "79-43="
$t.Cell(16,3).Range.Text = "77-65="
$t.Cell(16,4).Range.Text = "68-42="
$t.Cell(16,5).Range.Text = "16-0="
$t.Cell(17,1).Range.Text = "35+13="
$t.Cell(17,2).Range.Text = "24+4="
$t.Cell(17,3).Range.Text = "79-24="
$t.Cell(17,4).Range.Text = "95-1="
$t.Cell(17,5).Range.Text = "76-2="
$t.Cell(18,1).Range.Text = "12+37="
$t.Cell(18,2).Range.Text = "79-35="
$t.Cell(18,3).Range.Text = "81-10="
$t.Cell(18,4).Range.Text = "83-76="
$t.Cell(18,5).Range.Text = "79+14="
$t.Cell(19,1).Range.Text = "21-4="
$t.Cell(19,2).Range.Text = "45+44="
$t.Cell(19,3).Range.Text = "11+11="
$t.Cell(19,4).Range.Text = "6+13="
$t.Cell(19,5).Range.Text = "93-36="
$t.Cell(20,1).Range.Text = "94-91="
$t.Cell(20,2).Range.Text = "15+82="
$t.Cell(20,3).Range.Text = "29+28="
$t.Cell(20,4).Range.Text = "31+31="
$t.Cell(20,5).Range.Text = "21+11="
